$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.064.10'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.16%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '1.798.78'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.72%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = '  -0.31%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''307.62'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.83%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''1.002'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.32%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = '''0.4206'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.72%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = '''0.3592'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.94%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = '''0.07105'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.49%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = '''0.8436'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.08%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = '''20.17'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.41%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = '1.801.40'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -6.02%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = '''5.293'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.60%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '''6.365'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.84%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '''0.06761'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.14%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''1.006'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.06%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '''80.16'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.83%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = '''0.000008707'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -4.24%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = '  -0.56%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '''15.01'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.93%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '27.006.28'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.20%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''5.059'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.72%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''11.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '2.006.79'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -5.07%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = '''1.923'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.50%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''152.88'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.20%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = '''18.10'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.43%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = '''5.020'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -6.22%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = '''113.20'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.65%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = '''1.646'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -12.83%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = '''0.09008'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.81%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = '''0.7245'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -8.21%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = '''2.866'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.12%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = '''4.331'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -6.41%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = '''1.088'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -7.74%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = '  -0.27%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = '''1.082'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.74%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = '''0.05134'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.87%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = '''0.01903'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.26%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = '''0.1627'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.12%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = '''0.4962'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.61%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = '''2.610'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -8.03%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = '''8.049'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -7.23%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = '''5.911'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -13.03%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''105.03'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.90%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = '  -0.39%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = '''10.17'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.41%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = '  -4.10%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = '''0.4514'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -6.25%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = '  -4.41%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = '''1.708'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -8.39%  '
$ws.Range("E51").Style = "Normal"
